$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 313, shifting existing rows 313:330 down to 314:331.
$ws.Rows.Item(313).Insert()

# Populate the newly inserted row 313 with the new data record.
$ws.Cells.Item(313, 1).Value = 9
$ws.Cells.Item(313, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(313, 3).Value = "Metropolitana"
$ws.Cells.Item(313, 4).Value = 44931
$ws.Cells.Item(313, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(313, 5).Value = 13
$ws.Cells.Item(313, 6).Value = 100112030
$ws.Cells.Item(313, 7).Value = "Poroto granado"
$ws.Cells.Item(313, 8).Value = "Sin especificar"
$ws.Cells.Item(313, 9).Value = "Primera"
$ws.Cells.Item(313, 10).Value = 70
$ws.Cells.Item(313, 11).Value = 43000
$ws.Cells.Item(313, 12).Value = 45000
$ws.Cells.Item(313, 13).Value = 44000
$ws.Cells.Item(313, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(313, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(313, 16).Value = 1760
$ws.Cells.Item(313, 17).Value = 25
$ws.Cells.Item(313, 18).Value = "Hortaliza"
